$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormatLocal = "@"
$c.Value = "308.89"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormatLocal = "@"
$c.Value = "0.28%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormatLocal = "@"
$c.Value = "40.94"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormatLocal = "@"
$c.Value = "0.27%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormatLocal = "@"
$c.Value = "5.121"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormatLocal = "@"
$c.Value = "1.51%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormatLocal = "@"
$c.Value = "0.07617"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormatLocal = "@"
$c.Value = "-0.17%"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormatLocal = "@"
$c.Value = "-0.21%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormatLocal = "@"
$c.Value = "0.9078"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormatLocal = "@"
$c.Value = "-0.12%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormatLocal = "@"
$c.Value = "0.1272"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormatLocal = "@"
$c.Value = "24.53%"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormatLocal = "@"
$c.Value = "2.02%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormatLocal = "@"
$c.Value = "0.09017"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormatLocal = "@"
$c.Value = "-0.92%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormatLocal = "@"
$c.Value = "0.04294"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormatLocal = "@"
$c.Value = "-2.30%"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormatLocal = "@"
$c.Value = "-1.09%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormatLocal = "@"
$c.Value = "0.001258"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormatLocal = "@"
$c.Value = "2.23%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormatLocal = "@"
$c.Value = "0.005783"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormatLocal = "@"
$c.Value = "-0.14%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormatLocal = "@"
$c.Value = "3.351"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormatLocal = "@"
$c.Value = "-0.47%"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormatLocal = "@"
$c.Value = "4.279"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormatLocal = "@"
$c.Value = "0.93%"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormatLocal = "@"
$c.Value = "0.47%"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormatLocal = "@"
$c.Value = "6.924"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormatLocal = "@"
$c.Value = "2.60%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormatLocal = "@"
$c.Value = "0.1393"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormatLocal = "@"
$c.Value = "2.63%"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormatLocal = "@"
$c.Value = "-1.29%"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormatLocal = "@"
$c.Value = "0.04049"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormatLocal = "@"
$c.Value = "-2.53%"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormatLocal = "@"
$c.Value = "5.10%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormatLocal = "@"
$c.Value = "0.004045"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormatLocal = "@"
$c.Value = "-1.26%"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormatLocal = "@"
$c.Value = "-1.95%"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormatLocal = "@"
$c.Value = "24.70%"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormatLocal = "@"
$c.Value = "0.02418"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormatLocal = "@"
$c.Value = "-0.12%"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormatLocal = "@"
$c.Value = "0.05219"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormatLocal = "@"
$c.Value = "0.90%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormatLocal = "@"
$c.Value = "0.007840"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormatLocal = "@"
$c.Value = "0.96%"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormatLocal = "@"
$c.Value = "-0.80%"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormatLocal = "@"
$c.Value = "-3.92%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormatLocal = "@"
$c.Value = "0.001934"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormatLocal = "@"
$c.Value = "-0.56%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormatLocal = "@"
$c.Value = "0.007355"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormatLocal = "@"
$c.Value = "-8.36%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormatLocal = "@"
$c.Value = "0.3364"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormatLocal = "@"
$c.Value = "9.82%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormatLocal = "@"
$c.Value = "0.00006885"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormatLocal = "@"
$c.Value = "8.18%"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormatLocal = "@"
$c.Value = "0.32%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormatLocal = "@"
$c.Value = "0.1384"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormatLocal = "@"
$c.Value = "2,741.61%"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormatLocal = "@"
$c.Value = "-31.78%"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormatLocal = "@"
$c.Value = "0.32%"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormatLocal = "@"
$c.Value = "0.32%"
$c.Style = "Normal"
